$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B and C hold plain (non-numeric-looking) text, safe to set directly.
# Columns D and E hold numeric-looking / percentage text that Excel would otherwise
# auto-convert to a real number; force them to Text, assign, then clear the temporary
# number-format override so no residual style is left on the cell.
function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 2
Set-TextValue 'D2' '64.763.18'
Set-TextValue 'E2' '  +2.06%  '

# Row 3
Set-TextValue 'D3' '3.465.44'
Set-TextValue 'E3' '  +2.00%  '

# Row 4
Set-TextValue 'E4' '  +0.05%  '

# Row 5
Set-TextValue 'D5' '577.06'
Set-TextValue 'E5' '  -0.19%  '

# Row 6
Set-TextValue 'D6' '161.79'
Set-TextValue 'E6' '  +4.10%  '

# Row 7
Set-TextValue 'D7' '1.00'
Set-TextValue 'E7' '  -0.03%  '

# Row 8
Set-TextValue 'D8' '3.467.65'
Set-TextValue 'E8' '  +1.96%  '

# Row 9
Set-TextValue 'D9' '0.582'
Set-TextValue 'E9' '  +9.06%  '

# Row 10
Set-TextValue 'E10' '  -1.58%  '

# Row 11
Set-TextValue 'E11' '  +4.85%  '

# Row 12
Set-TextValue 'E12' '  +1.91%  '

# Row 13
Set-TextValue 'D13' '4.061.60'
Set-TextValue 'E13' '  +2.09%  '

# Row 14
Set-TextValue 'E14' '  -2.94%  '

# Row 15
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D15' '0.0000195'
Set-TextValue 'E15' '  +5.75%  '

# Row 16
$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D16' '28.99'
Set-TextValue 'E16' '  +7.18%  '

# Row 17
Set-TextValue 'D17' '64.754.51'
Set-TextValue 'E17' '  +1.90%  '

# Row 18
Set-TextValue 'D18' '3.475.68'
Set-TextValue 'E18' '  +3.62%  '

# Row 19
Set-TextValue 'E19' '  +0.39%  '

# Row 20
Set-TextValue 'D20' '14.46'
Set-TextValue 'E20' '  +3.29%  '

# Row 21
Set-TextValue 'D21' '393.04'
Set-TextValue 'E21' '  +1.41%  '

# Row 22
Set-TextValue 'E22' '  -2.52%  '

# Row 23
Set-TextValue 'D23' '0.550'
Set-TextValue 'E23' '  +2.61%  '

# Row 24
Set-TextValue 'D24' '73.01'
Set-TextValue 'E24' '  +3.36%  '

# Row 25
Set-TextValue 'E25' '  +0.37%  '

# Row 26
Set-TextValue 'D26' '0.0000125'
Set-TextValue 'E26' '  +20.79%  '

# Row 27
Set-TextValue 'D27' '9.49'
Set-TextValue 'E27' '  -1.11%  '

# Row 28
Set-TextValue 'D28' '0.181'
Set-TextValue 'E28' '  +0.65%  '

# Row 29
Set-TextValue 'E29' '  +0.12%  '

# Row 30
Set-TextValue 'D30' '6.18'
Set-TextValue 'E30' '  +9.93%  '

# Row 31
Set-TextValue 'E31' '  +8.13%  '

# Row 32
Set-TextValue 'E32' '  +0.06%  '

# Row 33
Set-TextValue 'D33' '6.58'

# Row 34
Set-TextValue 'D34' '23.72'
Set-TextValue 'E34' '  +2.67%  '

# Row 35
Set-TextValue 'D35' '0.999'
Set-TextValue 'E35' '  +0.09%  '

# Row 36
Set-TextValue 'D36' '7.09'
Set-TextValue 'E36' '  +5.80%  '

# Row 37
Set-TextValue 'E37' '  +1.28%  '

# Row 38
Set-TextValue 'D38' '161.71'
Set-TextValue 'E38' '  +2.21%  '

# Row 39
Set-TextValue 'E39' '  +1.18%  '

# Row 40
Set-TextValue 'D40' '0.0775'
Set-TextValue 'E40' '  +1.94%  '

# Row 41
Set-TextValue 'D41' '27.71'
Set-TextValue 'E41' '  +0.59%  '

# Row 42
Set-TextValue 'D42' '2.939.62'
Set-TextValue 'E42' '  +1.24%  '

# Row 43
Set-TextValue 'D43' '4.59'
Set-TextValue 'E43' '  +6.57%  '

# Row 44
Set-TextValue 'E44' '  -1.15%  '

# Row 45
Set-TextValue 'D45' '42.85'
Set-TextValue 'E45' '  +3.47%  '

# Row 46
Set-TextValue 'D46' '0.776'
Set-TextValue 'E46' '  +1.69%  '

# Row 47
Set-TextValue 'D47' '24.37'
Set-TextValue 'E47' '  +9.69%  '

# Row 48
Set-TextValue 'D48' '1.10'
Set-TextValue 'E48' '  +2.97%  '

# Row 49
Set-TextValue 'D49' '2.21'
Set-TextValue 'E49' '  +14.71%  '

# Row 50
Set-TextValue 'D50' '0.876'
Set-TextValue 'E50' '  +7.50%  '

# Row 51
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D51' '6.59'
Set-TextValue 'E51' '  +4.21%  '
